# TemperatureCalculation.xlsx update
# - Replace the UI/API Celsius temperature readings used by the variance
#   calculation (B2, B3) with new source values coming from the external
#   "Comparator" utility. These are supplied as single-precision floats,
#   matching the numeric precision produced by that tool, so downstream
#   formulas (x-u, (x-u)^2, Sum, Variance, VarianceLogic) recompute
#   automatically.
# - Leave the cursor/selection on cell E9, matching where the workbook
#   was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [single]22.0
$ws.Range("B3").Value = [single]25.58

$ws.Range("E9").Select()
